$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - INCONT L.A. 4MG 30 TAB.
$ws.Range("H14").Value = "1:0"
$ws.Range("Q14").Value = "0:2"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "83.1600"
$ws.Range("P14").NumberFormat = "0.00"

# Row 18 - TAMSULIN 0.4MG 28 CAPS
$ws.Range("H18").Value = "1:1"
$ws.Range("Q18").Value = "1:0"
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "124.0000"
$ws.Range("P18").NumberFormat = "0.00"

# Row 29 - total sale price
$ws.Range("P29").Value = 798.81

# Updated report generation timestamp
$ws.Range("A30").Value = "Tuesday, 29 July, 2025 10:33 AM"
